# Actualización SmartScore desde Streamlit (Gerarado Juan)
#
# 1) Row 2: the 9 SmartScore cells that were stored as quoted-text numbers
#    ("0.620", "0.541", ...) become real numeric values.
# 2) A brand-new row 3 is appended with the full survey submission for
#    "Gerarado Juan" (the SmartScore cells in the *new* row stay textual,
#    exactly like the original export format, only row 2 got "fixed").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 2 — convert the inline-string SmartScore numbers to real numbers
# ---------------------------------------------------------------------
$ws.Range("I2").Value  = 0.62
$ws.Range("L2").Value  = 0.541
$ws.Range("O2").Value  = 0.532
$ws.Range("R2").Value  = 0.632
$ws.Range("U2").Value  = 0.597
$ws.Range("X2").Value  = 0.554
$ws.Range("AA2").Value = 0.678
$ws.Range("AD2").Value = 0.472
$ws.Range("AG2").Value = 0.441

# ---------------------------------------------------------------------
# 2) Row 3 — new participant "Gerarado Juan"
# ---------------------------------------------------------------------

# The SmartScore columns for this new row must stay as *text* (mirrors the
# original exporter, which always wrote them as quoted strings). Mark those
# cells as Text before writing so Excel does not auto-coerce them to numbers.
$smartScoreCols = @("I3","L3","O3","R3","U3","X3","AA3","AD3","AG3")
foreach ($addr in $smartScoreCols) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A3").Value = "Gerarado Juan_20251113_181055"
# B3 (Grupo_Experimental) is blank for this participant, same as row 2.
$ws.Range("C3").Value = "Gerarado Juan"
$ws.Range("D3").Value = 21
$ws.Range("E3").Value = "Male"
$ws.Range("F3").Value = "2025-11-13 18:10:55"

$pesosJson = @"
{
  "portion": 0.4,
  "diet": 0.7142857142857143,
  "salt": 0.4,
  "fat": 0.8,
  "natural": 0.8,
  "convenience": 0.4,
  "price": 0.6
}
"@
$ws.Range("G3").Value = $pesosJson

$ws.Range("H3").Value  = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I3").Value  = "0.563"
$ws.Range("J3").Value  = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Range("K3").Value  = "Maruchan Ramen Sabor Pollo"
$ws.Range("L3").Value  = "0.454"
$ws.Range("M3").Value  = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"

$ws.Range("N3").Value  = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("O3").Value  = "0.418"
$ws.Range("P3").Value  = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

$ws.Range("Q3").Value  = "Kraft Macaroni & Cheese Dinner"
$ws.Range("R3").Value  = "0.663"
$ws.Range("S3").Value  = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Range("T3").Value  = "Annie’s Shells & White Cheddar"
$ws.Range("U3").Value  = "0.608"
$ws.Range("V3").Value  = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

$ws.Range("W3").Value  = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("X3").Value  = "0.576"
$ws.Range("Y3").Value  = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

$ws.Range("Z3").Value  = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA3").Value = "0.730"
$ws.Range("AB3").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

$ws.Range("AC3").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AD3").Value = "0.580"
$ws.Range("AE3").Value = "Portátil, saludable, fácil, buena textura, sabor suave"

$ws.Range("AF3").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AG3").Value = "0.556"
$ws.Range("AH3").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# Writing the multi-line "Pesos" JSON blob makes Excel auto-expand the row
# height to fit every embedded line break; re-running AutoFit snaps both
# data rows back to the sheet's normal (default) height.
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
